# New certificate about C Programming added
# (row 5: LP Academy / "C Programming For Beginners - Master the C Language")
# The course name + academy were already present; this fills in the
# completion date and the certificate link for that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Completion date for the C Programming certificate (2021-07-03),
# written as the underlying date serial so the cell's existing date
# number-format (already applied to column D) is preserved as-is.
$ws.Range("D5").Value = 44380

# Certificate link for the C Programming certificate
$certUrl = "https://www.udemy.com/certificate/UC-acb11d70-fc58-4488-b4a8-3b5161b07040/"
$ws.Range("E5").Value = $certUrl
$ws.Hyperlinks.Add($ws.Range("E5"), $certUrl, [Type]::Missing, [Type]::Missing, $certUrl)

# Move the active selection (as recorded in the saved file) to E6
$ws.Range("E6").Select()
